# Updates cryptos list price/volume cells (columns D and E, rows 2-51)
# per the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.091.20"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.764.87"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.70"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5237"
$ws.Range("E7").Value = "  +3.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2754"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("E9").Value = "  -3.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06202"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "1.776.27"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07010"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.74"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6412"
$ws.Range("E14").Value = "  +6.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.535"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.22"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "26.109.32"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006753"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("D22").Value = "2.002.96"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.075"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.442"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.195"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.90"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.508"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.850"
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.18"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.13"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08409"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.705"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.458"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04464"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.623"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6063"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.749"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.984"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.72"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3883"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7454"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.940"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05515"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.364"
$ws.Range("E47").Value = "  +6.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1121"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.23"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.72"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.74%  "
